$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Sát Phá Tham"
$ws.Range("B3").Value = "Là người chủ về sát phạt"
$ws.Range("C3").Value = "Đầu óc có tính thực tế và muốn hành động nhanh, ngay, dứt khoát."

$ws.Range("C3").Select()
